$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A47").Value = "L:"
$ws.Range("A48").Value = "16.99, 36.01, 68.54, 40.2, 54.79, 57.88, 42.47, 24.88, 2.47, 10.08, 34.11, 41.02, 6.64, 4.67, 5.87, 41.34, 161.88, 38.19, 69.9, 6.19, 40.93, 139.26, 2239.9, 843.3, 53.45, 97.99, 1582.14, 32.07, 33.6, 68.14, 33, 70.51, 640.75, 54.36, 47.21"
$ws.Range("A49").Value = "R:"
$ws.Range("A50").Value = "330.4, 154.07, 217.84, 105.83, 5182.37, 1.04, 808.16"
